$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("BasicSearch")
$ws2 = $wb.Worksheets.Item("AdvanceSearch")

# --- BasicSearch sheet: refactor the Age column into ValidAges/InvalidAges ---

# Insert a new column C (ZipCode column shifts from C -> D), inheriting
# formatting from the (old) C column to its left.
$ws1.Columns.Item(3).Insert()

# Header row
$ws1.Cells.Item(1, 2).Value = "ValidAges"
$ws1.Cells.Item(1, 3).Value = "InvalidAges"

# Column widths for the two age columns
$ws1.Columns.Item(2).ColumnWidth = 9.5
$ws1.Columns.Item(3).ColumnWidth = 14.6666666666667

# Row 2: Breast Cancer
$ws1.Cells.Item(2, 2).Formula = "'1"
$ws1.Cells.Item(2, 3).Formula = "'0"

# Row 3: Brain Cancer (ValidAges already holds 50 from old data)
$ws1.Cells.Item(3, 3).Formula = "'-3"

# Row 4: Lung Cancer (ValidAges already holds 120 from old data)
$ws1.Cells.Item(4, 3).Value = 121

# New row 5 with a non-numeric, invalid age value
$ws1.Cells.Item(5, 3).Value = "chicken"

# Update selection / active cell on BasicSearch and make it the active sheet
$ws1.Activate()
$ws1.Range("C3").Select()
